$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.615.32"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "3.001.93"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "3.001.52"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "3.496.34"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "61.623.11"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "3.002.04"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("D35").Value = "0.0₃0832"
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +10.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "402.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value = "2.704.86"
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("E51").Value = "  +2.10%  "
